$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update cryptocurrency price values in column D, preserving their
# original text (inline string) representation.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "242.84"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "23.07"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "5.411"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "3.437"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "6.526"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.8087"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.9390"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.1425"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07385"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.03287"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.03064"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.09346"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.844"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.001579"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.04665"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.0005909"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.005947"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.001255"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.004903"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.00006799"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "3.558"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.127"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.03970"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.006180"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1074"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.002570"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.009490"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.00005230"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.6698"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.002328"
